# Patch novendra, quiz sebelum perbaikan view
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Document")

# Update the image filename text values (shared strings) referenced by C8 and C20
$ws.Range("C8").Value = "be_assets\quiz\quiz_0f13a20c-16f0-4fbf-a33f-c6cecf985e19.png"
$ws.Range("C20").Value = "be_assets\quiz\quiz_9fe95eef-eb1c-4350-a67f-d437a7203f5a.png"

# Update the sheet view: scroll position and selection
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("F8").Select()
